# Generate Report for Handoff
# b.md has now been handed off for localization (zh-cn and de-de).
# Update the three worksheets (Overview, zh-cn, de-de) to reflect the
# new "Ready for handoff" status for the b.md row, including the new
# handoff file names, handoff datetimes, and (for the locale sheets)
# the "content duplicate" flag and error detail message.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet: row 3 is b.md
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2017-02-17 08:20:12"

# ---------------------------------------------------------------
# zh-cn sheet: row 3 is b.md
# ---------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2017-02-17 08:19:54"
$zhcn.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5efd40fd4f9bb8d0ad15359747991c5834034a25/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e0f48a10ed1af438d86bd1ad77ff46ce01b846ae/e2e/b.md."
$zhcn.Columns.Item(18).ColumnWidth = 40

# ---------------------------------------------------------------
# de-de sheet: row 3 is b.md
# ---------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2017-02-17 08:20:12"
$dede.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5efd40fd4f9bb8d0ad15359747991c5834034a25/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e0f48a10ed1af438d86bd1ad77ff46ce01b846ae/e2e/b.md."
$dede.Columns.Item(18).ColumnWidth = 40
